$d = $word.ActiveDocument
$range = $d.Content
$range.Find.Execute(" problem ", $false, $false, $false, $false, $false, $true, 1, $false, " defect ", 2)
